$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 747.5714
$ws.Range("I6").Value = 607.6667
$ws.Range("J6").Value = 785.7273
$ws.Range("K6").Value = 1823.0001
$ws.Range("L6").Value = 2357.1819
$ws.Range("M6").Value = -1711.0001
$ws.Range("N6").Value = -2581.1819

$ws.Range("H33").Value = 830.7
$ws.Range("I33").Value = 702.2
$ws.Range("J33").Value = 959.2
$ws.Range("K33").Value = 702.2
$ws.Range("L33").Value = 959.2
$ws.Range("M33").Value = -473.2
$ws.Range("N33").Value = -1417.2

$ws.Range("H76").Value = 4837.1113
$ws.Range("I76").Value = 4294.45
$ws.Range("K76").Value = 4294.45
$ws.Range("M76").Value = -3979.45

$ws.Range("H79").Value = 4837.1113
$ws.Range("I79").Value = 4294.45
$ws.Range("K79").Value = 4294.45
$ws.Range("M79").Value = -3202.45

$ws.Range("H98").Value = 1884.4375
$ws.Range("I98").Value = 1048.9231
$ws.Range("K98").Value = 1048.9231
$ws.Range("M98").Value = 449.0769

$ws.Range("H122").Value = 1884.4375
$ws.Range("I122").Value = 1048.9231
$ws.Range("K122").Value = 3146.7693
$ws.Range("M122").Value = -696.7692999999999

$ws.Range("H137").Value = 8207.98
$ws.Range("I137").Value = 11730.667
$ws.Range("K137").Value = 35192.001
$ws.Range("M137").Value = -32642.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14671.437
$ws.Range("I32").Value = 14825.439
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 14825.439
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -14538.439
$ws.Range("N32").Value = -10574

$ws.Range("H61").Value = 3058.4138
$ws.Range("I61").Value = 2112.818
$ws.Range("K61").Value = 2112.818
$ws.Range("M61").Value = -1900.818

$ws.Range("H74").Value = 195224.28
$ws.Range("I74").Value = 206625.73
$ws.Range("K74").Value = 206625.73
$ws.Range("M74").Value = -205751.73

$ws.Range("H77").Value = 195224.28
$ws.Range("I77").Value = 206625.73
$ws.Range("K77").Value = 1033128.65
$ws.Range("M77").Value = -1028760.65

$ws.Range("H132").Value = 1170.721
$ws.Range("I132").Value = 938.1539
$ws.Range("J132").Value = 3438.25
$ws.Range("K132").Value = 2814.4617
$ws.Range("L132").Value = 10314.75
$ws.Range("M132").Value = -284.4616999999998
$ws.Range("N132").Value = -15374.75

$ws.Range("H136").Value = 3058.4138
$ws.Range("I136").Value = 2112.818
$ws.Range("K136").Value = 6338.454000000001
$ws.Range("M136").Value = -3788.454000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H101").Value = 44577
$ws.Range("J101").Value = 44577
$ws.Range("L101").Value = 44577
$ws.Range("N101").Value = -51067

$ws.Range("H105").Value = 7203.8
$ws.Range("I105").Value = 520
$ws.Range("J105").Value = 8874.75
$ws.Range("K105").Value = 520
$ws.Range("L105").Value = 8874.75
$ws.Range("M105").Value = 1227
$ws.Range("N105").Value = -12368.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1679.8529
$ws.Range("I58").Value = 1467.6
$ws.Range("K58").Value = 1467.6
$ws.Range("M58").Value = -1264.6

$ws.Range("H62").Value = 12233.167
$ws.Range("I62").Value = 4600
$ws.Range("K62").Value = 4600
$ws.Range("M62").Value = -3976

$ws.Range("H65").Value = 12233.167
$ws.Range("I65").Value = 4600
$ws.Range("K65").Value = 23000
$ws.Range("M65").Value = -19880

$ws.Range("H109").Value = 39994
$ws.Range("J109").Value = 39994
$ws.Range("L109").Value = 39994
$ws.Range("N109").Value = -42074

$ws.Range("H136").Value = 1679.8529
$ws.Range("I136").Value = 1467.6
$ws.Range("K136").Value = 4402.799999999999
$ws.Range("M136").Value = -1852.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1438.5294
$ws.Range("I107").Value = 1771
$ws.Range("J107").Value = 1064.5
$ws.Range("K107").Value = 5313
$ws.Range("L107").Value = 3193.5
$ws.Range("M107").Value = -3393
$ws.Range("N107").Value = -7033.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 6024200
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 6024200
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 6024200
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -6025392

$ws.Range("H70").Value = 7202.4375
$ws.Range("I70").Value = 7039.3335
$ws.Range("J70").Value = 7412.143
$ws.Range("K70").Value = 7039.3335
$ws.Range("L70").Value = 7412.143
$ws.Range("M70").Value = -6769.3335
$ws.Range("N70").Value = -7952.143

$ws.Range("H73").Value = 7202.4375
$ws.Range("I73").Value = 7039.3335
$ws.Range("J73").Value = 7412.143
$ws.Range("K73").Value = 7039.3335
$ws.Range("L73").Value = 7412.143
$ws.Range("M73").Value = -6103.3335
$ws.Range("N73").Value = -9284.143

$ws.Range("H126").Value = 2418
$ws.Range("I126").Value = 2366.647
$ws.Range("J126").Value = 2542.7144
$ws.Range("K126").Value = 7099.941
$ws.Range("L126").Value = 7628.1432
$ws.Range("M126").Value = -4629.941
$ws.Range("N126").Value = -12568.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2004
$ws.Range("I7").Value = 2004
$ws.Range("K7").Value = 2004
$ws.Range("M7").Value = -1892

$ws.Range("H40").Value = 4021
$ws.Range("I40").Value = 4021
$ws.Range("K40").Value = 4021
$ws.Range("M40").Value = -3885

$ws.Range("H103").Value = 31500
$ws.Range("J103").Value = 31500
$ws.Range("L103").Value = 31500
$ws.Range("N103").Value = -33844

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H126").Value = 2004
$ws.Range("I126").Value = 2004
$ws.Range("K126").Value = 6012
$ws.Range("M126").Value = -3542

$ws.Range("H136").Value = 3065.625
$ws.Range("I136").Value = 2526.923
$ws.Range("K136").Value = 7580.768999999999
$ws.Range("M136").Value = -5030.768999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 87500
$ws.Range("J109").Value = 87500
$ws.Range("L109").Value = 87500
$ws.Range("N109").Value = -90274

$ws.Range("H113").Value = 1149.625
$ws.Range("I113").Value = 758.7
$ws.Range("K113").Value = 2276.1
$ws.Range("M113").Value = -106.1000000000004

$ws.Range("H122").Value = 43274
$ws.Range("I122").Value = 63954.477
$ws.Range("K122").Value = 191863.431
$ws.Range("M122").Value = -189413.431
